# Table S3 Deep sequencing statistics - restructure columns
#
# Summary of the edit:
#   - Swap column pair (B,C) [Dataset, Sample] with column pair (D,E)
#     [Library size, Antibody format] for rows 1-21 (full table incl. header),
#     preserving values, number formats, fonts, borders, alignment and merges.
#   - Rename header B1 "Library size" -> "Theoretical library size" and turn on
#     wrap text for that cell.
#   - Rename header F1 "Number of reads replicate 1" -> "NiR1" and
#     H1 "Number of reads replicate 2" -> "NiR2".
#   - Header row is no longer bold.
#   - Adjust column widths for B,C,D,E to match the new content.
#   - Update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the (B:C) and (D:E) column blocks (rows 1-21) ---------------
# Stage old D:E (Library size / Antibody format) in a scratch area far to
# the right, then move old B:C into D:E, then move the staged content into
# B:C.  Using Copy() (not Cut()) keeps formatting/merges intact and lets us
# fully control the clean-up of the scratch range afterwards.

$ws.Range("D1:E21").Copy($ws.Range("P1:Q21"))
$ws.Range("B1:C21").Copy($ws.Range("D1:E21"))
$ws.Range("P1:Q21").Copy($ws.Range("B1:C21"))
$ws.Range("P1:Q21").Clear()

# --- 2. Header text updates ------------------------------------------------

$ws.Range("B1").Value = "Theoretical library size"
$ws.Range("B1").WrapText = $true

$ws.Range("F1").Value = "NiR1"
$ws.Range("H1").Value = "NiR2"

# --- 3. Header row is no longer bold ---------------------------------------

$ws.Range("A1:J1").Font.Bold = $false

# --- 4. Column widths --------------------------------------------------

$ws.Columns(2).ColumnWidth = 14.666666666666666   # B width 15.5
$ws.Columns(3).ColumnWidth = 7.666666666666667    # C width 8.5 (no bestFit)
$ws.Columns(4).ColumnWidth = 7.666666666666667    # D width 8.5
$ws.Columns(5).ColumnWidth = 5.498697916666667    # E width 6.33203125

# --- 5. Restore the active cell selection -----------------------------

$ws.Range("K8").Select()
